# Apply the crypto price/volume update from the GitHub Actions scrape run.
# Values are written as text (NumberFormat "@") to stop Excel from
# auto-coercing numeric-looking strings (e.g. "307.53", "0.999") into
# real numbers, matching the original inline-string cell content; the
# format is then restored to "General" so cell formatting is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $newValue
    $range.NumberFormat = "General"
}

Set-TextValue "D2" "40.039.33"
Set-TextValue "E2" "  -3.93%  "
Set-TextValue "D3" "2.331.60"
Set-TextValue "E3" "  -5.63%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.07%  "
Set-TextValue "D5" "307.53"
Set-TextValue "E5" "  -3.98%  "
Set-TextValue "D6" "85.08"
Set-TextValue "E6" "  -7.48%  "
Set-TextValue "D7" "0.528"
Set-TextValue "E7" "  -3.80%  "
Set-TextValue "E8" "  +0.02%  "
Set-TextValue "D9" "0.484"
Set-TextValue "E9" "  -4.94%  "
Set-TextValue "D10" "0.0818"
Set-TextValue "E10" "  -3.77%  "
Set-TextValue "D11" "30.08"
Set-TextValue "E11" "  -8.54%  "
Set-TextValue "E12" "  +0.27%  "
Set-TextValue "D13" "2.692.95"
Set-TextValue "E13" "  -5.56%  "
Set-TextValue "D14" "6.42"
Set-TextValue "E14" "  -6.55%  "
Set-TextValue "D15" "14.72"
Set-TextValue "E15" "  -4.50%  "
Set-TextValue "D16" "2.337.94"
Set-TextValue "E16" "  -5.48%  "
Set-TextValue "D17" "0.753"
Set-TextValue "E17" "  -4.65%  "
Set-TextValue "D18" "40.011.42"
Set-TextValue "E18" "  -3.78%  "
Set-TextValue "D19" "0.0₃0903"
Set-TextValue "E19" "  -3.74%  "
Set-TextValue "D20" "6.08"
Set-TextValue "E20" "  -5.37%  "
Set-TextValue "D21" "67.60"
Set-TextValue "E21" "  -5.29%  "
Set-TextValue "D22" "10.68"
Set-TextValue "E22" "  -4.64%  "
Set-TextValue "D23" "235.24"
Set-TextValue "E23" "  -1.55%  "
Set-TextValue "E24" "  -6.87%  "
Set-TextValue "E25" "  +0.12%  "
Set-TextValue "D26" "1.80"
Set-TextValue "E26" "  -6.98%  "
Set-TextValue "D27" "23.40"
Set-TextValue "E27" "  -5.92%  "
Set-TextValue "E28" "  -1.20%  "
Set-TextValue "D29" "9.22"
Set-TextValue "E29" "  -5.05%  "
Set-TextValue "D30" "34.95"
Set-TextValue "E30" "  -3.78%  "
Set-TextValue "D31" "152.29"
Set-TextValue "E31" "  -2.88%  "
Set-TextValue "E32" "  -0.11%  "
Set-TextValue "D33" "5.13"
Set-TextValue "E33" "  -5.46%  "
Set-TextValue "B34" "WEMIXToken"
Set-TextValue "C34" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D34" "2.44"
Set-TextValue "E34" "  -4.62%  "
Set-TextValue "B35" "Hedera"
Set-TextValue "C35" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D35" "0.0724"
Set-TextValue "E35" "  -5.11%  "
Set-TextValue "E36" "  -1.92%  "
Set-TextValue "D37" "0.0998"
Set-TextValue "E37" "  -3.23%  "
Set-TextValue "B38" "LidoDAOToken"
Set-TextValue "C38" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D38" "2.75"
Set-TextValue "E38" "  -4.87%  "
Set-TextValue "B39" "Celestia"
Set-TextValue "C39" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D39" "15.74"
Set-TextValue "E39" "  -7.62%  "
Set-TextValue "D40" "1.70"
Set-TextValue "E40" "  -6.81%  "
Set-TextValue "D41" "3.81"
Set-TextValue "E41" "  -4.71%  "
Set-TextValue "E42" "  -5.70%  "
Set-TextValue "D43" "1.944.79"
Set-TextValue "E43" "  -2.81%  "
Set-TextValue "E44" "  -5.15%  "
Set-TextValue "D45" "17.63"
Set-TextValue "E45" "  -5.50%  "
Set-TextValue "D46" "9.26"
Set-TextValue "E46" "  -1.71%  "
Set-TextValue "D47" "2.68"
Set-TextValue "E47" "  -9.21%  "
Set-TextValue "D48" "2.564.80"
Set-TextValue "E48" "  -6.06%  "
Set-TextValue "D49" "92.81"
Set-TextValue "E49" "  -4.69%  "
Set-TextValue "D50" "71.35"
Set-TextValue "E50" "  -5.65%  "
Set-TextValue "E51" "  -2.41%  "
